$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44186
$ws.Range("L2").Value = 'Especial'
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 22500
$ws.Range("O2").Value = 23000
$ws.Range("P2").Value = 22750
$ws.Range("Q2").Value = '$/caja 18 kilos'
$ws.Range("S2").Value = 1264
$ws.Range("T2").Value = 18
$ws.Range("D3").Value = 44189
$ws.Range("L3").Value = 'Primera'
$ws.Range("M3").Value = 400
$ws.Range("N3").Value = 23500
$ws.Range("O3").Value = 24000
$ws.Range("P3").Value = 23750
$ws.Range("Q3").Value = '$/caja 18 kilos'
$ws.Range("R3").Value = 'Región de O''Higgins'
$ws.Range("S3").Value = 1319
$ws.Range("T3").Value = 18
$ws.Range("D4").Value = 44189
$ws.Range("K4").Value = 'Dina'
$ws.Range("L4").Value = 'Segunda'
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 21500
$ws.Range("O4").Value = 22000
$ws.Range("P4").Value = 21750
$ws.Range("Q4").Value = '$/caja 18 kilos'
$ws.Range("R4").Value = 'Región de O''Higgins'
$ws.Range("S4").Value = 1208
$ws.Range("T4").Value = 18
$ws.Range("D5").Value = 44537
$ws.Range("N5").Value = 22000
$ws.Range("O5").Value = 23000
$ws.Range("P5").Value = 22500
$ws.Range("Q5").Value = '$/caja 18 kilos'
$ws.Range("S5").Value = 1250
$ws.Range("T5").Value = 18
$ws.Range("D6").Value = 44537
$ws.Range("L6").Value = 'Segunda'
$ws.Range("M6").Value = 240
$ws.Range("N6").Value = 18000
$ws.Range("O6").Value = 19000
$ws.Range("P6").Value = 18500
$ws.Range("Q6").Value = '$/caja 18 kilos'
$ws.Range("S6").Value = 1028
$ws.Range("T6").Value = 18
$ws.Range("D7").Value = 44162
$ws.Range("L7").Value = 'Especial'
$ws.Range("N7").Value = 20500
$ws.Range("O7").Value = 21000
$ws.Range("P7").Value = 20750
$ws.Range("Q7").Value = '$/caja 15 kilos'
$ws.Range("S7").Value = 1383
$ws.Range("T7").Value = 15
$ws.Range("D8").Value = 44162
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 300
$ws.Range("N8").Value = 17500
$ws.Range("O8").Value = 18000
$ws.Range("P8").Value = 17750
$ws.Range("Q8").Value = '$/caja 15 kilos'
$ws.Range("S8").Value = 1183
$ws.Range("T8").Value = 15
$ws.Range("D9").Value = 44174
$ws.Range("K9").Value = 'Castle Brite'
$ws.Range("L9").Value = 'Primera'
$ws.Range("N9").Value = 22500
$ws.Range("P9").Value = 22750
$ws.Range("S9").Value = 1264
$ws.Range("D10").Value = 44181
$ws.Range("K10").Value = 'Modesto'
$ws.Range("M10").Value = 16
$ws.Range("N10").Value = 495000
$ws.Range("O10").Value = 500000
$ws.Range("P10").Value = 497500
$ws.Range("Q10").Value = '$/bins (500 kilos)'
$ws.Range("S10").Value = 995
$ws.Range("T10").Value = 500
$ws.Range("D11").Value = 44181
$ws.Range("K11").Value = 'Modesto'
$ws.Range("M11").Value = 10
$ws.Range("N11").Value = 425000
$ws.Range("O11").Value = 430000
$ws.Range("P11").Value = 427500
$ws.Range("Q11").Value = '$/bins (500 kilos)'
$ws.Range("S11").Value = 855
$ws.Range("T11").Value = 500
$ws.Range("D12").Value = 44165
$ws.Range("L12").Value = 'Especial'
$ws.Range("N12").Value = 20500
$ws.Range("O12").Value = 21000
$ws.Range("P12").Value = 20750
$ws.Range("Q12").Value = '$/caja 15 kilos'
$ws.Range("S12").Value = 1383
$ws.Range("T12").Value = 15
$ws.Range("D13").Value = 44165
$ws.Range("L13").Value = 'Primera'
$ws.Range("M13").Value = 200
$ws.Range("N13").Value = 17500
$ws.Range("O13").Value = 18000
$ws.Range("P13").Value = 17750
$ws.Range("R13").Value = 'Región Metropolitana'
$ws.Range("S13").Value = 1183
$ws.Range("D14").Value = 44539
$ws.Range("L14").Value = 'Especial'
$ws.Range("M14").Value = 160
$ws.Range("N14").Value = 24500
$ws.Range("O14").Value = 25000
$ws.Range("P14").Value = 24750
$ws.Range("S14").Value = 1650
$ws.Range("D15").Value = 44539
$ws.Range("L15").Value = 'Primera'
$ws.Range("M15").Value = 160
$ws.Range("N15").Value = 22500
$ws.Range("O15").Value = 23000
$ws.Range("P15").Value = 22750
$ws.Range("S15").Value = 1517
$ws.Range("D16").Value = 44539
$ws.Range("K16").Value = 'Castle Brite'
$ws.Range("L16").Value = 'Segunda'
$ws.Range("N16").Value = 18000
$ws.Range("O16").Value = 18500
$ws.Range("P16").Value = 18250
$ws.Range("Q16").Value = '$/caja 15 kilos'
$ws.Range("R16").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S16").Value = 1217
$ws.Range("T16").Value = 15
$ws.Range("D17").Value = 44532
$ws.Range("L17").Value = 'Primera'
$ws.Range("M17").Value = 200
$ws.Range("N17").Value = 21000
$ws.Range("O17").Value = 22000
$ws.Range("P17").Value = 21500
$ws.Range("R17").Value = 'Región Metropolitana'
$ws.Range("S17").Value = 1433
$ws.Range("D18").Value = 44168
$ws.Range("M18").Value = 200
$ws.Range("N18").Value = 23500
$ws.Range("O18").Value = 24000
$ws.Range("P18").Value = 23750
$ws.Range("Q18").Value = '$/caja 18 kilos'
$ws.Range("R18").Value = 'Región Metropolitana'
$ws.Range("S18").Value = 1319
$ws.Range("T18").Value = 18
$ws.Range("D19").Value = 44167
$ws.Range("L19").Value = 'Especial'
$ws.Range("M19").Value = 400
$ws.Range("N19").Value = 20000
$ws.Range("O19").Value = 21000
$ws.Range("P19").Value = 20500
$ws.Range("R19").Value = 'Región Metropolitana'
$ws.Range("S19").Value = 1367
$ws.Range("D20").Value = 44167
$ws.Range("L20").Value = 'Primera'
$ws.Range("M20").Value = 360
$ws.Range("N20").Value = 17000
$ws.Range("O20").Value = 18000
$ws.Range("P20").Value = 17500
$ws.Range("S20").Value = 1167
$ws.Range("D21").Value = 44161
$ws.Range("K21").Value = 'Dina'
$ws.Range("M21").Value = 240
$ws.Range("N21").Value = 19500
$ws.Range("O21").Value = 20000
$ws.Range("P21").Value = 19750
$ws.Range("S21").Value = 1317
$ws.Range("D22").Value = 44161
$ws.Range("K22").Value = 'Dina'
$ws.Range("L22").Value = 'Segunda'
$ws.Range("M22").Value = 140
$ws.Range("N22").Value = 17500
$ws.Range("O22").Value = 18000
$ws.Range("P22").Value = 17750
$ws.Range("Q22").Value = '$/caja 15 kilos'
$ws.Range("S22").Value = 1183
$ws.Range("T22").Value = 15
$ws.Range("D23").Value = 44187
$ws.Range("K23").Value = 'Dina'
$ws.Range("M23").Value = 240
$ws.Range("N23").Value = 22000
$ws.Range("O23").Value = 23000
$ws.Range("P23").Value = 22500
$ws.Range("Q23").Value = '$/caja 18 kilos'
$ws.Range("S23").Value = 1250
$ws.Range("T23").Value = 18
$ws.Range("D24").Value = 44536
$ws.Range("M24").Value = 160
$ws.Range("N24").Value = 22000
$ws.Range("O24").Value = 23000
$ws.Range("P24").Value = 22500
$ws.Range("Q24").Value = '$/caja 18 kilos'
$ws.Range("S24").Value = 1250
$ws.Range("T24").Value = 18
$ws.Range("D25").Value = 44536
$ws.Range("K25").Value = 'Castle Brite'
$ws.Range("L25").Value = 'Segunda'
$ws.Range("M25").Value = 120
$ws.Range("N25").Value = 18000
$ws.Range("O25").Value = 19000
$ws.Range("P25").Value = 18500
$ws.Range("Q25").Value = '$/caja 18 kilos'
$ws.Range("S25").Value = 1028
$ws.Range("T25").Value = 18
$ws.Range("D26").Value = 44540
$ws.Range("K26").Value = 'Castle Brite'
$ws.Range("L26").Value = 'Especial'
$ws.Range("M26").Value = 140
$ws.Range("N26").Value = 24500
$ws.Range("O26").Value = 25000
$ws.Range("P26").Value = 24750
$ws.Range("Q26").Value = '$/caja 15 kilos'
$ws.Range("R26").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S26").Value = 1650
$ws.Range("T26").Value = 15
$ws.Range("D27").Value = 44540
$ws.Range("K27").Value = 'Castle Brite'
$ws.Range("M27").Value = 100
$ws.Range("N27").Value = 22500
$ws.Range("O27").Value = 23000
$ws.Range("P27").Value = 22750
$ws.Range("Q27").Value = '$/caja 15 kilos'
$ws.Range("R27").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S27").Value = 1517
$ws.Range("T27").Value = 15
$ws.Range("D28").Value = 44540
$ws.Range("K28").Value = 'Castle Brite'
$ws.Range("M28").Value = 100
$ws.Range("N28").Value = 18000
$ws.Range("O28").Value = 18500
$ws.Range("P28").Value = 18250
$ws.Range("Q28").Value = '$/caja 15 kilos'
$ws.Range("R28").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S28").Value = 1217
$ws.Range("T28").Value = 15
